{"js": "// The diff splits the paragraph that ends with the sentence\n//   \"...ef \u00fe\u00fa ert \u201eh\u00e6ttur a\u00f0 vinna\"\"\n// right before the trailing <w:bookmarkStart _GoBack/><w:bookmarkEnd/>\n// pair: that paragraph now ends right after its last run, and a brand\n// new (otherwise empty) centered paragraph - carrying the same\n// sz/szCs=24 run-properties - is inserted after it to hold the\n// \"_GoBack\" bookmark on its own.\n\nconst body = context.document.body;\n\n// Locate the existing \"_GoBack\" bookmark range before touching anything;\n// this is the exact split point the diff targets.\nconst bookmarkRange = body.getBookmarkRange(\"_GoBack\");\n\n// Remove the old bookmark - it will be re-created on the new paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// OOXML for the new paragraph: centered, sz/szCs 24 (matching the\n// paragraph it is split off from), containing only the bookmark.\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr>' +\n  '<w:jc w:val=\"center\"/>' +\n  '<w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"1\"/>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\n// Insert the new paragraph right after the (now bookmark-less) point\n// where the bookmark used to live - this splits the original paragraph\n// in two, exactly as the diff shows.\nbookmarkRange.insertOoxml(newParagraphOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The diff splits the paragraph that ends with the sentence\n#   \"...ef \u00fe\u00fa ert \u201eh\u00e6ttur a\u00f0 vinna\"\"\n# right before the trailing bookmarkStart/bookmarkEnd pair for\n# \"_GoBack\": that paragraph now ends right after its last run, and a\n# brand new (otherwise empty) centered paragraph - carrying the same\n# sz/szCs=24 run-properties - is inserted after it to hold the\n# \"_GoBack\" bookmark on its own.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that currently owns the \"_GoBack\" bookmark -\n# this is the paragraph we need to split right at its end.\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$targetParagraph = $bm.Range.Paragraphs.Item(1)\n\n# Split the paragraph in two: everything stays in the first paragraph,\n# a new (empty) paragraph is created right after it. At this point the\n# \"_GoBack\" bookmark is still attached to the end of the first\n# (original) paragraph.\n$targetParagraph.Range.InsertParagraphAfter()\n\n# Re-fetch the paragraph collection: the new, empty paragraph is the\n# one right after the original paragraph.\n$newParaIndex = $targetParagraph.Range.Paragraphs.Item(1).Index + 1\n$newParagraph = $d.Paragraphs.Item($newParaIndex)\n\n# The new paragraph is completely empty (just its paragraph mark), and\n# adding a bookmark exactly at that lone boundary position is\n# unreliable, so temporarily add a placeholder character to give the\n# bookmark a real, unambiguous insertion point.\n$newParagraph.Range.InsertAfter(\"x\")\n$newParagraph = $d.Paragraphs.Item($newParaIndex)\n$bookmarkPoint = $d.Range($newParagraph.Range.Start, $newParagraph.Range.Start)\n\n# Re-adding a bookmark with an existing name moves it here (Word\n# replaces same-named bookmarks), so the old one on the first\n# paragraph goes away automatically.\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n\n# Remove the placeholder character again, leaving the new paragraph\n# with nothing but its (now relocated) bookmark.\n$newParagraph = $d.Paragraphs.Item($newParaIndex)\n$placeholderRange = $d.Range($newParagraph.Range.Start, $newParagraph.Range.Start + 1)\n$placeholderRange.Delete()\n"}
